$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.170.02"
$ws.Range("E2").Value = "  +1.29%  "

$ws.Range("D3").Value = "2.276.43"
$ws.Range("E3").Value = "  +0.45%  "

$ws.Range("E4").Value = "  -0.22%  "

$ws.Range("D5").Value = "231.28"
$ws.Range("E5").Value = "  -0.35%  "

$ws.Range("D6").Value = "0.635"
$ws.Range("E6").Value = "  +1.70%  "

$ws.Range("D7").Value = "63.97"
$ws.Range("E7").Value = "  +4.56%  "

$ws.Range("E8").Value = "  +0.01%  "

$ws.Range("D9").Value = "0.449"
$ws.Range("E9").Value = "  +9.12%  "

$ws.Range("E10").Value = "  +9.96%  "

$ws.Range("D11").Value = "57.16"
$ws.Range("E11").Value = "  -0.80%  "

$ws.Range("D12").Value = "26.91"
$ws.Range("E12").Value = "  +19.08%  "

$ws.Range("E13").Value = "  +1.53%  "

$ws.Range("D14").Value = "2.616.81"
$ws.Range("E14").Value = "  +0.55%  "

$ws.Range("D15").Value = "15.71"
$ws.Range("E15").Value = "  -0.01%  "

$ws.Range("E16").Value = "  +7.64%  "

$ws.Range("D17").Value = "0.840"
$ws.Range("E17").Value = "  +3.82%  "

$ws.Range("D18").Value = "2.277.34"
$ws.Range("E18").Value = "  +0.55%  "

$ws.Range("D19").Value = "44.018.25"
$ws.Range("E19").Value = "  +1.32%  "

$ws.Range("D20").Value = "0.0000101"
$ws.Range("E20").Value = "  +7.78%  "

$ws.Range("D21").Value = "73.82"
$ws.Range("E21").Value = "  +1.09%  "

$ws.Range("E22").Value = "  -1.45%  "

$ws.Range("D23").Value = "252.92"
$ws.Range("E23").Value = "  +1.82%  "

$ws.Range("E24").Value = "  +0.14%  "

$ws.Range("D25").Value = "2.45"
$ws.Range("E25").Value = "  -5.39%  "

$ws.Range("B26").Value = "Toncoin"
$ws.Range("C26").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D26").Value = "2.31"
$ws.Range("E26").Value = "  -3.52%  "

$ws.Range("B27").Value = "WEMIXToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D27").Value = "3.36"
$ws.Range("E27").Value = "  +24.82%  "

$ws.Range("D28").Value = "10.03"
$ws.Range("E28").Value = "  +2.35%  "

$ws.Range("D29").Value = "171.96"
$ws.Range("E29").Value = "  +1.35%  "

$ws.Range("E30").Value = "  -2.82%  "

$ws.Range("D31").Value = "20.84"
$ws.Range("E31").Value = "  +1.17%  "

$ws.Range("E32").Value = "  -6.81%  "

$ws.Range("E33").Value = "  +2.99%  "

$ws.Range("D34").Value = "0.0699"
$ws.Range("E34").Value = "  +6.65%  "

$ws.Range("D35").Value = "4.82"
$ws.Range("E35").Value = "  +1.66%  "

$ws.Range("E36").Value = "  -2.47%  "

$ws.Range("E37").Value = "  +5.31%  "

$ws.Range("D38").Value = "6.56"
$ws.Range("E38").Value = "  +1.47%  "

$ws.Range("E39").Value = "  -3.15%  "

$ws.Range("E40").Value = "  +3.58%  "

$ws.Range("E41").Value = "  +8.40%  "

$ws.Range("E42").Value = "  -0.31%  "

$ws.Range("D43").Value = "17.66"
$ws.Range("E43").Value = "  +5.32%  "

$ws.Range("D44").Value = "0.0975"
$ws.Range("E44").Value = "  +0.12%  "

$ws.Range("D45").Value = "8.25"
$ws.Range("E45").Value = "  -5.17%  "

$ws.Range("D46").Value = "10.41"
$ws.Range("E46").Value = "  +17.85%  "

$ws.Range("D47").Value = "98.37"
$ws.Range("E47").Value = "  +1.12%  "

$ws.Range("D48").Value = "1.21"
$ws.Range("E48").Value = "  -0.63%  "

$ws.Range("D49").Value = "4.38"
$ws.Range("E49").Value = "  -0.40%  "

$ws.Range("D50").Value = "1.449.24"
$ws.Range("E50").Value = "  -1.30%  "

$ws.Range("E51").Value = "  +2.94%  "
